$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: R2,R5 resistor switched from 1206 12.4K part to 0805 part ---
$ws.Range("B7").Value = "RES 12.4K OHM 1% 1/4W 0805"
$ws.Range("C7").Value = "RNCP0805FTD12K4"

# --- Row 8: R3,R6 resistor switched from 1206 9.09K part to 0805 part ---
$ws.Range("B8").Value = "RES 9.09K OHM 1% 1/4W 0805"
$ws.Range("C8").Value = "RNCP0805FTD9K09"

# --- Row 9: R7 resistor switched from 1206 158K part to 0805 part, new datasheet link ---
$ws.Range("B9").Value = "RES SMD 158K OHM 0.5% 1/4W 0805"
$ws.Range("C9").Value = "  ERJ-PB6D1583V "
$ws.Hyperlinks.Add($ws.Range("D9"), "https://industrial.panasonic.com/cdbs/www-data/pdf/RDM0000/AOA0000C328.pdf")
$ws.Range("D9").Value = "https://industrial.panasonic.com/cdbs/www-data/pdf/RDM0000/AOA0000C328.pdf"

# --- Row 12: C3 cap switched from 1206 1nF part to 0805 1000pF part, new datasheet link ---
$ws.Range("B12").Value = "CAP CER 1000PF 25V X7R 0805"
$ws.Range("C12").Value = "VJ0805Y102JXXPW1BC"
$ws.Hyperlinks.Add($ws.Range("D12"), "https://www.vishay.com/docs/28548/vjw1bcbascomseries.pdf")
$ws.Range("D12").Value = "https://www.vishay.com/docs/28548/vjw1bcbascomseries.pdf"

# --- Row 11: C1,C2,C4,C5,C6 cap switched from 1206 0.1uF part to 0805 part, new datasheet link ---
$ws.Range("D11").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D11"), "http://datasheets.avx.com/AutoMLCC.pdf")
$ws.Range("D11").Value = "http://datasheets.avx.com/AutoMLCC.pdf"
$ws.Range("B11").Value = "CAP CER 0.1UF 50V X7R 0805"
$ws.Range("C11").Value = "08055C104J4T2A"

# --- Row 10: R8,R9 resistor switched from 1206 2.32K part to 0805 part ---
$ws.Range("B10").Value = "RES 2.32K OHM 1% 1/4W 0805"
$ws.Range("C10").Value = "RNCP0805FTD2K32"

# --- Datasheet links for rows 7 and 10 now point at SEI-rncp instead of yageo ---
$ws.Range("D7").Value = "https://www.seielect.com/Catalog/SEI-rncp.pdf"
$ws.Hyperlinks.Add($ws.Range("D10"), "https://www.seielect.com/Catalog/SEI-rncp.pdf")
$ws.Range("D10").Value = "https://www.seielect.com/Catalog/SEI-rncp.pdf"

# --- Restore the plain (non-bold-applyFont) hyperlink-look style on the cells whose
#     style slot got touched by Hyperlinks.Add, so they reuse the existing style
#     instead of cloning a near-duplicate one ---
$ws.Range("D6").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D12").PasteSpecial(-4122)

# --- Updated unit prices ---
$ws.Range("E7").Value = 0.1
$ws.Range("E8").Value = 0.1
$ws.Range("E9").Value = 0.25
$ws.Range("E10").Value = 0.1
$ws.Range("E11").Value = 0.27
$ws.Range("E12").Value = 0.34

# --- Rows 9 & 10 now carry an explicit (default) row height ---
$ws.Rows.Item(9).RowHeight = 14.4
$ws.Rows.Item(10).RowHeight = 14.4

# --- Update the active selection left by the editor ---
$null = $ws.Range("B19").Select()
